# Integrate Contact page Test Scenario into TestData.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two unused sheets (Sheet2, Sheet3) ---
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# --- The pre-existing A5/D5 cells lose their old "special" styling and become ---
# --- plain bordered cells, like the rest of the body rows                    ---
$ws.Range("A5").Font.Bold = $false
$ws.Range("D5").Font.Bold = $false

# --- New scenario row: TC05 create contact with mandatory fields          ---
# --- (written in the same left-to-right order the row was authored in)    ---
$ws.Range("D6").Value = "nature"
$ws.Range("I1").Value = "ContactName"
$ws.Range("I6").Value = "Test_Enginner"
$ws.Range("E6").Value = "Sachin"
$ws.Range("F6").Value = "Tendulkar"
$ws.Range("A6").Value = "TC05_create_contact_with_mandatory_fields"
$ws.Range("B6").Value = "admin"
$ws.Range("C6").Value = "admin"
$ws.Range("H6").Value = "ACC"

# --- New "ContactName" header in column I: same bold+fill header look as the ---
# --- rest of row 1, plus a full border on all four sides                     ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Borders.LineStyle = 1
$ws.Columns.Item(9).ColumnWidth = 12.22

# --- Give the rest of column I (and the new G6) the same plain bordered look ---
# --- used by the other empty body cells                                     ---
$ws.Range("D2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("G6").PasteSpecial(-4122)

# --- Update the active selection like the saved workbook shows ---
$ws.Range("C20").Select()
